$wb = $excel.ActiveWorkbook

# --- Reuse existing sheet, rename to "Events" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Events"

# Remove the bits of the old layout (11 cols x 2 rows) that fall outside
# the new layout (8 cols x 6 rows) so no stale data/formatting lingers.
$ws1.Range("I1:K2").Clear()
$ws1.Range("F2:H2").ClearContents()

# --- Add a new sheet "Stats" right after "Events" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Stats"

# =========================================================
# Sheet 1: Events
# =========================================================

# Header row
$ws1.Range("A1").Value = "Datum"
$ws1.Range("B1").Value = "Schul/Uni Name"
$ws1.Range("C1").Value = "Adresse"
$ws1.Range("D1").Value = "Stadt"
$ws1.Range("E1").Value = "Bundesland"
$ws1.Range("F1").Value = "PLZ"
$ws1.Range("G1").Value = "Tische"
$ws1.Range("H1").Value = "Teilnehmer"

# Row 2
$ws1.Range("A2").Value = "23.01.2024"
$ws1.Range("B2").Value = "Hochschule Karlsruhe"
$ws1.Range("C2").Value = "Moltkestraße 30"
$ws1.Range("D2").Value = "Karlsruhe"
$ws1.Range("E2").Value = "Baden-Württemberg"
$ws1.Range("F2").Value = 76133
$ws1.Range("G2").Value = 3
$ws1.Range("H2").Value = 15

# Row 3
$ws1.Range("A3").Value = "23.01.2024"
$ws1.Range("B3").Value = "Universität Stuttgart"
$ws1.Range("C3").Value = "Arminstraße 15"
$ws1.Range("D3").Value = "Stuttgart"
$ws1.Range("E3").Value = "Baden-Württemberg"
$ws1.Range("F3").Value = 70178
$ws1.Range("G3").Value = 10
$ws1.Range("H3").Value = 50

# Row 4
$ws1.Range("A4").Value = "23.01.2024"
$ws1.Range("B4").Value = "Hochschule München"
$ws1.Range("C4").Value = "Lothstraße 34"
$ws1.Range("D4").Value = "München"
$ws1.Range("E4").Value = "Bayern"
$ws1.Range("F4").Value = 80335
$ws1.Range("G4").Value = 5
$ws1.Range("H4").Value = 25

# Row 5
$ws1.Range("A5").Value = "23.01.2024"
$ws1.Range("B5").Value = "PH Karlsruhe"
$ws1.Range("C5").Value = "Moltkestraße 10"
$ws1.Range("D5").Value = "Karlsruhe"
$ws1.Range("E5").Value = "Baden-Württemberg"
$ws1.Range("F5").Value = 76133
$ws1.Range("G5").Value = 4
$ws1.Range("H5").Value = 16

# Row 6 (PLZ stored as text here, matching source data)
$ws1.Range("A6").Value = "23.01.2024"
$ws1.Range("B6").Value = "Hochschule Karlsruhe"
$ws1.Range("C6").Value = "Moltkestraße 30"
$ws1.Range("D6").Value = "Karlsruhe"
$ws1.Range("E6").Value = "Baden-Württemberg"
$ws1.Range("F6").NumberFormat = "@"
$ws1.Range("F6").Value = "76133"
$ws1.Range("G6").Value = 2
$ws1.Range("H6").Value = 12

# =========================================================
# Sheet 2: Stats
# =========================================================

# Header row
$ws2.Range("A1").Value = "Schul/Uni Name"
$ws2.Range("B1").Value = "Stadt"
$ws2.Range("C1").Value = "PLZ"
$ws2.Range("D1").Value = "Latitude"
$ws2.Range("E1").Value = "Longitude"
$ws2.Range("F1").Value = "EventCount"
$ws2.Range("G1").Value = "CityEventTotal"
$ws2.Range("H1").Value = "TotalTables"
$ws2.Range("I1").Value = "TotalParticipants"

# Row 2
$ws2.Range("A2").Value = "Hochschule Karlsruhe"
$ws2.Range("B2").Value = "Karlsruhe"
$ws2.Range("C2").Value = 76133
$ws2.Range("D2").Value = 49.013238
$ws2.Range("E2").Value = 8.392054
$ws2.Range("F2").Value = 2
$ws2.Range("G2").Value = 3
$ws2.Range("H2").Value = 5
$ws2.Range("I2").Value = 27

# Row 3
$ws2.Range("A3").Value = "Universität Stuttgart"
$ws2.Range("B3").Value = "Stuttgart"
$ws2.Range("C3").Value = 70178
$ws2.Range("D3").Value = 48.7647049
$ws2.Range("E3").Value = 9.166575080084748
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 10
$ws2.Range("I3").Value = 50

# Row 4
$ws2.Range("A4").Value = "Hochschule München"
$ws2.Range("B4").Value = "München"
$ws2.Range("C4").Value = 80335
$ws2.Range("D4").Value = 48.154141
$ws2.Range("E4").Value = 11.55624838157751
$ws2.Range("F4").Value = 1
$ws2.Range("G4").Value = 1
$ws2.Range("H4").Value = 5
$ws2.Range("I4").Value = 25

# Row 5
$ws2.Range("A5").Value = "PH Karlsruhe"
$ws2.Range("B5").Value = "Karlsruhe"
$ws2.Range("C5").Value = 76133
$ws2.Range("D5").Value = 49.01473905
$ws2.Range("E5").Value = 8.3948733
$ws2.Range("F5").Value = 1
$ws2.Range("G5").Value = 3
$ws2.Range("H5").Value = 4
$ws2.Range("I5").Value = 16

# =========================================================
# Header styling: bold, centered/top aligned, thin box border
# (build the style once on Events!A1:H1, then clone it via
#  copy/paste-special so both sheets share the same style record)
# =========================================================

$headerRange1 = $ws1.Range("A1:H1")
$headerRange1.Font.Bold = $true
$headerRange1.HorizontalAlignment = -4108
$headerRange1.VerticalAlignment = -4160
$headerRange1.Borders.LineStyle = 1

$ws1.Range("A1").Copy()
$ws2.Range("A1:I1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$null = $ws1.Activate()
$null = $ws1.Range("A1").Select()
